$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '35.579.99'
$ws.Range("E2").Value = '  +0.93%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.899.05'
$ws.Range("E3").Value = '  -0.52%  '
$ws.Range("E4").Value = '  -0.80%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.69'
$ws.Range("E5").Value = '  -2.95%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.694'
$ws.Range("E6").Value = '  -3.96%  '
$ws.Range("E7").Value = '  -0.84%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '43.97'
$ws.Range("E8").Value = '  +8.35%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.354'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0744'
$ws.Range("E10").Value = '  -2.13%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0969'
$ws.Range("E11").Value = '  -2.01%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '13.18'
$ws.Range("E12").Value = '  +2.64%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.173.65'
$ws.Range("E13").Value = '  -0.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.729'
$ws.Range("E14").Value = '  +0.21%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.95'
$ws.Range("E15").Value = '  -0.33%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.918.09'
$ws.Range("E16").Value = '  +0.56%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '35.489.40'
$ws.Range("E17").Value = '  +0.67%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '73.94'
$ws.Range("E18").Value = '  -1.26%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0827'
$ws.Range("E19").Value = '  -2.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '248.02'
$ws.Range("E20").Value = '  +1.78%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.94'
$ws.Range("E21").Value = '  -0.82%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.99'
$ws.Range("E22").Value = '  -2.25%  '
$ws.Range("E23").Value = '  -0.72%  '
$ws.Range("E24").Value = '  +3.97%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.21'
$ws.Range("E25").Value = '  -9.36%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.75'
$ws.Range("E26").Value = '  -0.28%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.52'
$ws.Range("E27").Value = '  -2.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.43'
$ws.Range("E28").Value = '  -1.62%  '
$ws.Range("E29").Value = '  -3.40%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.128.47'
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.79'
$ws.Range("E31").Value = '  +9.09%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.27'
$ws.Range("E32").Value = '  -2.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0582'
$ws.Range("E33").Value = '  -0.99%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.24'
$ws.Range("E34").Value = '  +0.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.862'
$ws.Range("E36").Value = '  -5.84%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.03'
$ws.Range("E37").Value = '  -0.56%  '
$ws.Range("E38").Value = '  -20.93%  '
$ws.Range("B39").Value = 'Aave'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '98.34'
$ws.Range("E39").Value = '  +1.99%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0682'
$ws.Range("E40").Value = '  +4.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.20'
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("E42").Value = '  -1.79%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.10'
$ws.Range("E43").Value = '  -2.18%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.39'
$ws.Range("E44").Value = '  -1.63%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.298.49'
$ws.Range("E45").Value = '  -2.82%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0802'
$ws.Range("E46").Value = '  +6.32%  '
$ws.Range("E47").Value = '  -1.17%  '
$ws.Range("E48").Value = '  -0.58%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '12.16'
$ws.Range("E49").Value = '  +3.87%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '43.66'
$ws.Range("E50").Value = '  -2.87%  '
$ws.Range("E51").Value = '  -5.32%  '
